$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Elimina los periodos de mora anteriores y agrega los nuevos: la tabla
# de periodos (columna E, filas 16-24) se reordena de mas reciente (2307)
# a mas antiguo (2211). El valor de mora (columna F) de 34666 queda
# asociado al periodo 2307 (ahora en la fila 16) y el periodo 2211 (ahora
# en la fila 24) toma el valor 40000.

$periodos = @("2307", "2306", "2305", "2304", "2303", "2302", "2301", "2212", "2211")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}

$ws.Cells.Item(16, 6).Value = 34666
$ws.Cells.Item(24, 6).Value = 40000
